$d = $word.ActiveDocument

# Change 1: merge the "lacks the ability to ... clientele." text into a single run
# (removes the grammar-check proofErr markers that wrapped
# "easily and seamlessly grow its clientele")
$d.Content.Find.Execute(
    "lacks the ability to easily and seamlessly grow its clientele. Mr",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "lacks the ability to easily and seamlessly grow its clientele. Mr",
    2)

# Change 2: remove the "with solutions provided by April 15, 2022" clause,
# leaving just "... improve."
$d.Content.Find.Execute(
    " improve, with solutions provided by April 15, 2022.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " improve.",
    2)
